$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Insert 10 new rows before the old "UC1 + config" block (old rows 13-22
# become rows 23-32's predecessors; the old row 17 block shifts to row 27).
$ws.Rows("13:22").Insert()

# ---- Row 13: UC2 / Create Classes ----
$ws.Cells.Item(13, 1).Value = "UC2"
$ws.Range("B13:B16").Clear()
$ws.Cells.Item(13, 2).Value = "Create Classes"
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(13, 4).Value = 6
$ws.Cells.Item(13, 5).Formula = "=(D13-C13)/C13"
$ws.Cells.Item(13, 6).Value = 6
$ws.Cells.Item(13, 7).Formula = "=(F13-D13)/D13"

# ---- Row 14: Create Properties ----
$ws.Cells.Item(14, 2).Value = "Create Properties"
$ws.Cells.Item(14, 3).Value = 30
$ws.Cells.Item(14, 4).Value = 34
$ws.Cells.Item(14, 5).Formula = "=(D14-C14)/C14"
$ws.Cells.Item(14, 6).Value = 30
$ws.Cells.Item(14, 7).Formula = "=(F14-D14)/D14"

# ---- Row 15: Create Associations ----
$ws.Cells.Item(15, 2).Value = "Create Associations"
$ws.Cells.Item(15, 3).Value = 24
$ws.Cells.Item(15, 4).Value = 24
$ws.Cells.Item(15, 5).Formula = "=(D15-C15)/C15"
$ws.Cells.Item(15, 6).Value = 9
$ws.Cells.Item(15, 7).Formula = "=(F15-D15)/D15"

# ---- Row 16: Create Operations ----
$ws.Cells.Item(16, 2).Value = "Create Operations"
$ws.Cells.Item(16, 3).Value = 34
$ws.Cells.Item(16, 4).Value = 34
$ws.Cells.Item(16, 5).Formula = "=(D16-C16)/C16"
$ws.Cells.Item(16, 6).Value = 28
$ws.Cells.Item(16, 7).Formula = "=(F16-D16)/D16"

# ---- Row 17: Sum of the UC2 block ----
$ws.Cells.Item(17, 3).Formula = "=SUM(C13:C16)"
$ws.Cells.Item(17, 4).Formula = "=SUM(D13:D16)"
$ws.Cells.Item(17, 5).Formula = "=(D17-C17)/C17"
$ws.Cells.Item(17, 6).Formula = "=SUM(F13:F16)"
$ws.Cells.Item(17, 7).Formula = "=(F17-D17)/D17"

# ---- Row 18: stray error row (C18/D18/F18 empty -> #DIV/0!) ----
$ws.Cells.Item(18, 5).Formula = "=(D18-C18)/C18"
$ws.Cells.Item(18, 7).Formula = "=(F18-D18)/D18"

# ---- Row 19-20: only G has the stray error formula ----
$ws.Cells.Item(19, 7).Formula = "=(F19-D19)/D19"
$ws.Cells.Item(20, 7).Formula = "=(F20-D20)/D20"

# ---- Rows 19-22: drop the leftover styled-but-empty cells so the rows
# match the target exactly (no content, no formula) ----
$ws.Range("E19:E22").Clear()
$ws.Range("G21:G22").Clear()

# Restore selection/view to match the target state
$ws.Range("C16").Select()
